$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-table updates: cell reference -> new text value.
# Using NumberFormat '@' + Value + Style 'Normal' to force text storage
# (matches source data which stores these as literal text, not parsed numbers/percentages),
# while avoiding leaving a quote-prefix/number-format style on the cell.
$updates = @(
    @('D2', '28.248.70'),
    @('E2', '  +1.27%  '),
    @('D3', '1.869.19'),
    @('E3', '  +3.35%  '),
    @('E4', '  +0.15%  '),
    @('D5', '312.03'),
    @('E5', '  +0.98%  '),
    @('D6', '1.002'),
    @('E6', '  +0.24%  '),
    @('D7', '0.5045'),
    @('E7', '  +1.13%  '),
    @('D8', '0.3914'),
    @('E8', '  +0.67%  '),
    @('D9', '0.09614'),
    @('E9', '  +1.14%  '),
    @('E10', '  +4.16%  '),
    @('D11', '40.81'),
    @('E11', '  +1.26%  '),
    @('D12', '6.466'),
    @('E12', '  +1.18%  '),
    @('D13', '21.01'),
    @('E13', '  +2.56%  '),
    @('D14', '1.867.45'),
    @('E14', '  +3.08%  '),
    @('E15', '  +0.13%  '),
    @('D16', '7.405'),
    @('E16', '  +1.99%  '),
    @('E17', '  +0.22%  '),
    @('D18', '92.87'),
    @('E18', '  -0.42%  '),
    @('E19', '  +0.35%  '),
    @('D20', '17.68'),
    @('E20', '  +3.08%  '),
    @('D21', '1.002'),
    @('E21', '  +0.25%  '),
    @('D22', '6.201'),
    @('E22', '  +4.32%  '),
    @('D23', '28.291.71'),
    @('E23', '  +1.23%  '),
    @('D24', '11.30'),
    @('E24', '  +1.24%  '),
    @('E25', '  +2.06%  '),
    @('D26', '2.579'),
    @('E26', '  +7.57%  '),
    @('D27', '2.086.83'),
    @('E27', '  +3.28%  '),
    @('D28', '21.23'),
    @('E28', '  +2.01%  '),
    @('D29', '158.84'),
    @('E29', '  +1.09%  '),
    @('D30', '127.65'),
    @('E30', '  -0.80%  '),
    @('D31', '0.1063'),
    @('E31', '  -0.92%  '),
    @('D32', '1.067'),
    @('E32', '  +1.21%  '),
    @('E33', '  -0.06%  '),
    @('D34', '3.621'),
    @('E34', '  -0.06%  '),
    @('D35', '0.06761'),
    @('E35', '  -0.81%  '),
    @('D36', '9.533'),
    @('E36', '  +6.20%  '),
    @('D37', '0.02418'),
    @('E37', '  +4.55%  '),
    @('E38', '  +2.02%  '),
    @('D39', '11.52'),
    @('E39', '  +1.27%  '),
    @('B40', 'TheSandbox'),
    @('C40', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'),
    @('D40', '0.6366'),
    @('E40', '  +1.73%  '),
    @('B41', 'InternetComputer(DFINITY)'),
    @('C41', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'),
    @('D41', '5.008'),
    @('E41', '  +1.13%  '),
    @('D42', '1.184'),
    @('E43', '  +0.27%  '),
    @('E44', '  +2.36%  '),
    @('D45', '0.5992'),
    @('E45', '  +1.55%  '),
    @('E46', '  -1.00%  '),
    @('D47', '3.662'),
    @('E47', '  -0.59%  '),
    @('D48', '2.006'),
    @('E48', '  +2.70%  '),
    @('D49', '123.59'),
    @('E49', '  -0.43%  '),
    @('D50', '1.199'),
    @('E50', '  +1.67%  '),
    @('D51', '0.06855'),
    @('E51', '  +1.35%  '),
)

foreach ($update in $updates) {
    $ref = $update[0]
    $val = $update[1]
    $rng = $ws.Range($ref)
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = 'Normal'
}

'Applied ' + $updates.Count + ' cell updates'
